# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages are now in sync with en-US: status strings,
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns are populated, and the affected columns are widened to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a7aad33685df75d37f7531af785c3ef277db43c/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status cells for both rows move from
# "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de columns so the longer status text fits.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# Latest Target File (I) - new hyperlink to a.md, styled like the other
# hyperlink cells in the sheet (underlined, custom blue).
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $aMdUrl, "", "", "a.md")
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $aMdUrl, "", "", "a.md")
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = 15570276

# Latest Handback File (J) - zh-cn xlf handback file name.
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Latest Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-08-25 12:38:48"
$wsZhCn.Range("K3").Value = "2016-08-25 12:38:48"

# Widen Status (C) and Latest Handback File (J) columns.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Latest Target File (I) - new hyperlink to a.md.
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $aMdUrl, "", "", "a.md")
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $aMdUrl, "", "", "a.md")
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = 15570276

# Latest Handback File (J) - de-de xlf handback file name.
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# Latest Handback DateTime (K) - de-de handback happened a bit later.
$wsDeDe.Range("K2").Value = "2016-08-25 12:38:54"
$wsDeDe.Range("K3").Value = "2016-08-25 12:38:54"

# Widen Status (C) and Latest Handback File (J) columns.
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
